# feat: add 2022-Q3 data
#
# 1. Insert a brand-new worksheet "2022-Q3" right after "总计" (before "2022-Q2"),
#    populated with the 2022-Q3 fund holdings data.
# 2. Insert the 2022-Q3 summary row into the "总计" sheet (right after the header,
#    before the existing 2022-Q2 row), pushing every other row down by one.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$zongji = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row values
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Reuse the bold/centered/bordered header style already present on "总计"!B1:D1
$zongji.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q3Rows = @(
    @{A=0; B="006122"; C="华安低碳生活混合A";                     D="3.48"; E="89.43"; F="3.60"; G="0.1253"; H=6},
    @{A=1; B="011144"; C="华安汇宏精选混合A";                     D="0.64"; E="89.51"; F="3.84"; G="0.0246"; H=4},
    @{A=2; B="011145"; C="华安汇宏精选混合C";                     D="0.19"; E="89.51"; F="3.84"; G="0.0073"; H=4},
    @{A=3; B="004321"; C="前海开源沪港深强国产业灵活配置混合";     D="0.11"; E="78.52"; F="5.17"; G="0.0057"; H=10},
    @{A=4; B="006477"; C="中邮沪港深精选混合";                     D="0.06"; E="90.21"; F="7.28"; G="0.0044"; H=4},
    @{A=5; B="014970"; C="华安低碳生活混合C";                     D="0.02"; E="89.43"; F="3.60"; G="0.0007"; H=6}
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row.A
    $q3.Cells.Item($r, 2).Value = "'" + $row.B
    $q3.Cells.Item($r, 3).Value = $row.C
    $q3.Cells.Item($r, 4).Value = "'" + $row.D
    $q3.Cells.Item($r, 5).Value = "'" + $row.E
    $q3.Cells.Item($r, 6).Value = "'" + $row.F
    $q3.Cells.Item($r, 7).Value = "'" + $row.G
    $q3.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Column A (the running index) uses the same centered/bordered style as the
# rest of the workbook's index columns - copy it down from "总计"!A2.
$zongji.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Step 2: insert the new 2022-Q3 summary row into "总计", right after the header
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

# The freshly-inserted row 2 doesn't carry the right per-column formatting
# (Insert() only copies the row directly above, and the header row has no
# value/style in column A). Re-apply formats explicitly:
#   - A2 needs the bordered/centered "index" style used by A3:A8
#   - B2:D2 need to be *plain* (no inherited style bleeding in from Insert())
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial($xlPasteFormats)
$zongji.Range("A3:D3").Copy()
$zongji.Range("B2:D2").PasteSpecial($xlPasteFormats)

$summaryRows = @(
    @{A=0; B="2022-Q3"; C=6;  D=0.17},
    @{A=1; B="2022-Q2"; C=10; D=0.41},
    @{A=2; B="2021-Q4"; C=14; D=5.06},
    @{A=3; B="2021-Q3"; C=11; D=2.11},
    @{A=4; B="2021-Q2"; C=12; D=4.55},
    @{A=5; B="2021-Q1"; C=18; D=6.7},
    @{A=6; B="2020-Q4"; C=18; D=9.220000000000001}
)

$r = 2
foreach ($row in $summaryRows) {
    $zongji.Cells.Item($r, 1).Value = $row.A
    $zongji.Cells.Item($r, 2).Value = $row.B
    $zongji.Cells.Item($r, 3).Value = $row.C
    $zongji.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}
